$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Four new test cases were written, increasing the Total and Automated
# Test Cases counts for the "Read" row from 6 to 10.
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 10

# Update the active cell selection to E3.
$ws.Range("E3").Select()
